$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The stats rows for Chris Gayle got reshuffled in the source data; swap the
# runs/balls/fours/sixes (C:F) between the row pairs below (row 7 is
# untouched) using Copy so the text-stored-as-number cell type is preserved.
$pairs = @(
    @(2, 4),
    @(3, 6),
    @(5, 8)
)

$scratch = $ws.Range("H1:K1")

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $row1 = $ws.Range("C${r1}:F${r1}")
    $row2 = $ws.Range("C${r2}:F${r2}")

    $row1.Copy($scratch)
    $row2.Copy($row1)
    $scratch.Copy($row2)
}

$scratch.ClearContents()
